$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.45
$ws.Range("C3").Value = -13.101
$ws.Range("E3").Value = 16.228
$ws.Range("C4").Value = -12.73
$ws.Range("E9").Value = 17.03
$ws.Range("A11").Value = -21.822
$ws.Range("A12").Value = -21.652
$ws.Range("C14").Value = -12.083
$ws.Range("A15").Value = -21.957
$ws.Range("E15").Value = 15.825
$ws.Range("E19").Value = 16.239
$ws.Range("E20").Value = 16.387
$ws.Range("E25").Value = 17.063
$ws.Range("C26").Value = -11.973
$ws.Range("A27").Value = -21.862
$ws.Range("E27").Value = 16.337
$ws.Range("A28").Value = -21.885
$ws.Range("E28").Value = 16.897
$ws.Range("E30").Value = 16.113
$ws.Range("A31").Value = -22.001
$ws.Range("C31").Value = -13.531
$ws.Range("A32").Value = -21.858
$ws.Range("E32").Value = 16.893
$ws.Range("C35").Value = -13.176
$ws.Range("A36").Value = -20.275
$ws.Range("C37").Value = -13.129
$ws.Range("A38").Value = -19.696
$ws.Range("C39").Value = -13.041
$ws.Range("C40").Value = -12.51
$ws.Range("E44").Value = 16.482
$ws.Range("C45").Value = -12.485
$ws.Range("A46").Value = -21.938
$ws.Range("E47").Value = 16.553
$ws.Range("C52").Value = -11.286
$ws.Range("A54").Value = -22.15
$ws.Range("A55").Value = -22.088
$ws.Range("A56").Value = -21.997
$ws.Range("C57").Value = -13.502
$ws.Range("E58").Value = 16.66
$ws.Range("E62").Value = 16.593
$ws.Range("A67").Value = -21.554
$ws.Range("A69").Value = -21.656
$ws.Range("A72").Value = -21.431
$ws.Range("A73").Value = -19.994
$ws.Range("E77").Value = 16.651
$ws.Range("E78").Value = 16.394
$ws.Range("C81").Value = -13.031
$ws.Range("A83").Value = -21.722
$ws.Range("C83").Value = -13.026
$ws.Range("E84").Value = 16.581
$ws.Range("A86").Value = -22.257
$ws.Range("E89").Value = 17.085
$ws.Range("A91").Value = -21.508
$ws.Range("E91").Value = 17.39
$ws.Range("E92").Value = 17.293
$ws.Range("A93").Value = -21.421
$ws.Range("E96").Value = 16.44600000000001
$ws.Range("A99").Value = -19.861
$ws.Range("C100").Value = -12.522
$ws.Range("C102").Value = -13.018
$ws.Range("E102").Value = 16.418
